$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Weekly Quantity": shift values in rows 4-19 up by 3 (dropping the
# old rows 4-6 values), then delete the now-duplicate trailing rows 20-23.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

$weekly = @(
  @(45158.99999999999, 1),
  @(45200.99999999999, 22),
  @(45221.99999999999, 23),
  @(45487.99999999999, 26),
  @(45494.99999999999, 2),
  @(45501.99999999999, 4),
  @(45508.99999999999, 10),
  @(45515.99999999999, 72),
  @(45522.99999999999, 8),
  @(45529.99999999999, 18),
  @(45536.99999999999, 28),
  @(45543.99999999999, 8),
  @(45550.99999999999, 4),
  @(45557.99999999999, 8),
  @(45564.99999999999, 30),
  @(45585.99999999999, 20)
)

$row = 4
foreach ($pair in $weekly) {
  $ws1.Cells.Item($row, 1).Value = $pair[0]
  $ws1.Cells.Item($row, 2).Value = $pair[1]
  $row = $row + 1
}

# Remove the old rows 20-23 which are no longer part of the data set.
$ws1.Range("A20:B23").EntireRow.Delete()

# ---------------------------------------------------------------------------
# Sheet "Monthly Trend": update a handful of values and shift rows 6-9 up
# by one, then delete the now-duplicate trailing row 10.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Monthly Trend")

$ws2.Cells.Item(3, 2).Value = 11

$monthly = @(
  @(45504.99999999999, 42),
  @(45535.99999999999, 126),
  @(45565.99999999999, 50),
  @(45596.99999999999, 20)
)

$row = 6
foreach ($pair in $monthly) {
  $ws2.Cells.Item($row, 1).Value = $pair[0]
  $ws2.Cells.Item($row, 2).Value = $pair[1]
  $row = $row + 1
}

$ws2.Range("A10:B10").EntireRow.Delete()
